$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name labels in column A (reflects reordering of the underlying
#     shared-string table so each row shows the correct country again) ---
$ws.Range("A1").Value = 'Datos actualizados a 23 de Junio de 2020 a las 22:24'
$ws.Range("A21").Value = 'Sudafrica'
$ws.Range("A22").Value = 'Canada'
$ws.Range("A65").Value = 'Camerun'
$ws.Range("A66").Value = 'Argelia'
$ws.Range("A101").Value = 'Costa Rica'
$ws.Range("A102").Value = 'Croacia'
$ws.Range("A103").Value = 'Cuba'
$ws.Range("A205").Value = 'Lesoto'
$ws.Range("A206").Value = 'San Cristobal y Nieves'
$ws.Range("A207").Value = 'Islas Turcas y Caicos'
$ws.Range("A208").Value = 'Islas Malvinas'
$ws.Range("A209").Value = 'Groenlandia'

# --- Update statistic values (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 2413739
$ws.Range("C4").Value = 25586
$ws.Range("D4").Value = 1006138
$ws.Range("E4").Value = 1284278
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 713
$ws.Range("H4").Value = 123323
$ws.Range("B7").Value = 456062
$ws.Range("C7").Value = 15612
$ws.Range("D7").Value = 258523
$ws.Range("E7").Value = 183056
$ws.Range("F7").Value = 0
$ws.Range("B21").Value = 106108
$ws.Range("C21").Value = 4518
$ws.Range("D21").Value = 55045
$ws.Range("E21").Value = 48961
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 111
$ws.Range("H21").Value = 2102
$ws.Range("B22").Value = 101905
$ws.Range("C22").Value = 268
$ws.Range("D22").Value = 64629
$ws.Range("E22").Value = 28823
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 17
$ws.Range("H22").Value = 8453
$ws.Range("B29").Value = 58141
$ws.Range("C29").Value = 1332
$ws.Range("D29").Value = 15535
$ws.Range("E29").Value = 40241
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 87
$ws.Range("H29").Value = 2365
$ws.Range("B30").Value = 51643
$ws.Range("C30").Value = 1003
$ws.Range("D30").Value = 24991
$ws.Range("E30").Value = 22378
$ws.Range("F30").Value = 0
$ws.Range("B51").Value = 21512
$ws.Range("C51").Value = 430
$ws.Range("D51").Value = 15869
$ws.Range("E51").Value = 5335
$ws.Range("F51").Value = 0
$ws.Range("B65").Value = 12270
$ws.Range("C65").Value = 229
$ws.Range("D65").Value = 7774
$ws.Range("E65").Value = 4183
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 5
$ws.Range("H65").Value = 313
$ws.Range("B66").Value = 12076
$ws.Range("C66").Value = 156
$ws.Range("D66").Value = 8674
$ws.Range("E66").Value = 2541
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 9
$ws.Range("H66").Value = 861
$ws.Range("B101").Value = 2368
$ws.Range("C101").Value = 91
$ws.Range("D101").Value = 1129
$ws.Range("E101").Value = 1227
$ws.Range("F101").Value = 0
$ws.Range("H101").Value = 12
$ws.Range("B102").Value = 2366
$ws.Range("C102").Value = 30
$ws.Range("D102").Value = 2142
$ws.Range("E102").Value = 117
$ws.Range("F102").Value = 0
$ws.Range("H102").Value = 107
$ws.Range("B103").Value = 2318
$ws.Range("C103").Value = 3
$ws.Range("D103").Value = 2123
$ws.Range("E103").Value = 110
$ws.Range("F103").Value = 0
$ws.Range("H103").Value = 85
$ws.Range("B120").Value = 1477
$ws.Range("C120").Value = 47
$ws.Range("D120").Value = 1213
$ws.Range("E120").Value = 246
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 7
$ws.Range("H120").Value = 18
$ws.Range("B164").Value = 231
$ws.Range("C164").Value = 12
$ws.Range("D164").Value = 94
$ws.Range("E164").Value = 130
$ws.Range("F164").Value = 0
$ws.Range("B205").Value = 17
$ws.Range("C205").Value = 5
$ws.Range("D205").Value = 2
$ws.Range("E205").Value = 15
$ws.Range("F205").Value = 0
$ws.Range("B206").Value = 15
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 15
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("H206").Value = 0
$ws.Range("B207").Value = 14
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 11
$ws.Range("E207").Value = 2
$ws.Range("F207").Value = 0
$ws.Range("H207").Value = 1
$ws.Range("B209").Value = 13
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 13
$ws.Range("E209").Value = 10
$ws.Range("F209").Value = 0
